$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5289
$ws1.Range("F6").Value = 801
$ws1.Range("F7").Value = 300

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 40

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5289
$ws4.Range("F6").Value = 801
$ws4.Range("F7").Value = 40
$ws4.Range("F8").Value = 300
